# Applies the author's update:
#  1. Column C ("Förändrad") date serial changes from 45184 to 45186 for every
#     data row (rows 2-138).
#  2. Every HYPERLINK() formula in columns S-Y (rows 2-13) gets a second
#     argument added: the "friendly name" text, which equals the row's
#     "Beteckning" value in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldDate = 45184
$newDate = 45186

for ($row = 1; $row -le $lastRow; $row++) {
    $cCell = $ws.Cells.Item($row, 3)   # Column C
    $cVal = $cCell.Value()
    if ($cVal -eq $null) {
        continue
    }

    $serial = $null
    if ($cVal -is [DateTime]) {
        $serial = $cVal.ToOADate()
    } else {
        $serial = [double]$cVal
    }

    if ($serial -eq $oldDate) {
        $cCell.Value = $newDate
    }
}

$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($row = 1; $row -le $lastRow; $row++) {
    $label = $ws.Cells.Item($row, 1).Value()   # Column A - Beteckning
    if ([string]::IsNullOrEmpty($label)) {
        continue
    }

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($row, $col)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $newFormula = '=HYPERLINK("' + $url + '", "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
